$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 124
$ws.Range("E2").Value = 145
$ws.Range("F2").Value = 85.51724137931035
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 145
$ws.Range("F3").Value = 13.79310344827586
$ws.Range("E4").Value = 145
$ws.Range("F4").Value = 0.6896551724137931
$ws.Range("E5").Value = 145
$ws.Range("E6").Value = 145
$ws.Range("E7").Value = 145
$ws.Range("E8").Value = 145
$ws.Range("E9").Value = 145
$ws.Range("E10").Value = 145
$ws.Range("E11").Value = 145
$ws.Range("E12").Value = 145
$ws.Range("E13").Value = 145
$ws.Range("D14").Value = 136
$ws.Range("E14").Value = 143
$ws.Range("F14").Value = 95.1048951048951
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = 143
$ws.Range("F15").Value = 4.195804195804196
$ws.Range("A16").Value = 13
$ws.Range("C16").Value = 'Synonyms'
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 143
$ws.Range("F16").Value = 0.6993006993006993
$ws.Range("A17").Value = 12
$ws.Range("C17").Value = 'Ethics & Morals'
$ws.Range("E17").Value = 143
$ws.Range("E18").Value = 143
$ws.Range("E19").Value = 143
$ws.Range("E20").Value = 143
$ws.Range("E21").Value = 143
$ws.Range("E22").Value = 143
$ws.Range("E23").Value = 143
$ws.Range("E24").Value = 143
$ws.Range("E25").Value = 143
$ws.Range("D26").Value = 56
$ws.Range("E26").Value = 71
$ws.Range("F26").Value = 78.87323943661971
$ws.Range("D27").Value = 14
$ws.Range("E27").Value = 71
$ws.Range("F27").Value = 19.71830985915493
$ws.Range("A28").Value = 30
$ws.Range("C28").Value = 'Direct Air Capture'
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 71
$ws.Range("F28").Value = 1.408450704225352
$ws.Range("A29").Value = 24
$ws.Range("C29").Value = 'Synonyms'
$ws.Range("E29").Value = 71
$ws.Range("A30").Value = 26
$ws.Range("C30").Value = 'Ocean fertilisation'
$ws.Range("E30").Value = 71
$ws.Range("A31").Value = 27
$ws.Range("C31").Value = 'Ocean Alkalinisation'
$ws.Range("E31").Value = 71
$ws.Range("A32").Value = 28
$ws.Range("C32").Value = 'Ethics & Morals'
$ws.Range("E32").Value = 71
$ws.Range("A33").Value = 29
$ws.Range("C33").Value = 'Enhanced Weathering'
$ws.Range("E33").Value = 71
$ws.Range("E34").Value = 71
$ws.Range("E35").Value = 71
$ws.Range("E36").Value = 71
$ws.Range("E37").Value = 71
$ws.Range("D38").Value = 50
$ws.Range("E38").Value = 109
$ws.Range("F38").Value = 45.87155963302752
$ws.Range("D39").Value = 23
$ws.Range("E39").Value = 109
$ws.Range("F39").Value = 21.10091743119266
$ws.Range("A40").Value = 40
$ws.Range("C40").Value = 'Enhanced Weathering'
$ws.Range("D40").Value = 11
$ws.Range("E40").Value = 109
$ws.Range("F40").Value = 10.09174311926606
$ws.Range("A41").Value = 45
$ws.Range("C41").Value = 'BECCS'
$ws.Range("D41").Value = 10
$ws.Range("E41").Value = 109
$ws.Range("F41").Value = 9.174311926605505
$ws.Range("A42").Value = 36
$ws.Range("C42").Value = 'Ethics & Morals'
$ws.Range("D42").Value = 4
$ws.Range("E42").Value = 109
$ws.Range("F42").Value = 3.669724770642202
$ws.Range("A43").Value = 38
$ws.Range("C43").Value = 'Ocean fertilisation'
$ws.Range("E43").Value = 109
$ws.Range("F43").Value = 2.752293577981652
$ws.Range("A44").Value = 44
$ws.Range("C44").Value = 'Biochar'
$ws.Range("E44").Value = 109
$ws.Range("F44").Value = 2.752293577981652
$ws.Range("A45").Value = 46
$ws.Range("C45").Value = 'Afforestation/reforestation'
$ws.Range("D45").Value = 3
$ws.Range("E45").Value = 109
$ws.Range("F45").Value = 2.752293577981652
$ws.Range("A46").Value = 39
$ws.Range("C46").Value = 'Ocean Alkalinisation'
$ws.Range("D46").Value = 2
$ws.Range("E46").Value = 109
$ws.Range("F46").Value = 1.834862385321101
$ws.Range("E47").Value = 109
$ws.Range("E48").Value = 109
$ws.Range("E49").Value = 109
$ws.Range("D50").Value = 73
$ws.Range("E50").Value = 88
$ws.Range("F50").Value = 82.95454545454545
$ws.Range("D51").Value = 10
$ws.Range("E51").Value = 88
$ws.Range("F51").Value = 11.36363636363636
$ws.Range("A52").Value = 52
$ws.Range("C52").Value = 'Ethics & Morals'
$ws.Range("E52").Value = 88
$ws.Range("F52").Value = 2.272727272727273
$ws.Range("A53").Value = 59
$ws.Range("C53").Value = 'Afforestation/reforestation'
$ws.Range("D53").Value = 2
$ws.Range("E53").Value = 88
$ws.Range("F53").Value = 2.272727272727273
$ws.Range("A54").Value = 51
$ws.Range("C54").Value = 'Ocean Alkalinisation'
$ws.Range("D54").Value = 1
$ws.Range("E54").Value = 88
$ws.Range("F54").Value = 1.136363636363636
$ws.Range("A55").Value = 49
$ws.Range("C55").Value = 'Soil Carbon Sequestration'
$ws.Range("E55").Value = 88
$ws.Range("E56").Value = 88
$ws.Range("E57").Value = 88
$ws.Range("E58").Value = 88
$ws.Range("E59").Value = 88
$ws.Range("E60").Value = 88
$ws.Range("E61").Value = 88
$ws.Range("D62").Value = 21
$ws.Range("E62").Value = 44
$ws.Range("F62").Value = 47.72727272727273
$ws.Range("D63").Value = 16
$ws.Range("E63").Value = 44
$ws.Range("F63").Value = 36.36363636363637
$ws.Range("D64").Value = 3
$ws.Range("E64").Value = 44
$ws.Range("F64").Value = 6.818181818181817
$ws.Range("E65").Value = 44
$ws.Range("F65").Value = 4.545454545454546
$ws.Range("E66").Value = 44
$ws.Range("F66").Value = 4.545454545454546
$ws.Range("E67").Value = 44
$ws.Range("E68").Value = 44
$ws.Range("E69").Value = 44
$ws.Range("E70").Value = 44
$ws.Range("E71").Value = 44
$ws.Range("E72").Value = 44
$ws.Range("E73").Value = 44
$ws.Range("D74").Value = 18
$ws.Range("E74").Value = 24
$ws.Range("F74").Value = 75
$ws.Range("E75").Value = 24
$ws.Range("F75").Value = 16.66666666666666
$ws.Range("E76").Value = 24
$ws.Range("F76").Value = 4.166666666666666
$ws.Range("E77").Value = 24
$ws.Range("F77").Value = 4.166666666666666
$ws.Range("E78").Value = 24
$ws.Range("E79").Value = 24
$ws.Range("E80").Value = 24
$ws.Range("E81").Value = 24
$ws.Range("E82").Value = 24
$ws.Range("E83").Value = 24
$ws.Range("E84").Value = 24
$ws.Range("E85").Value = 24
$ws.Range("D86").Value = 3
$ws.Range("E86").Value = 3
$ws.Range("E87").Value = 3
$ws.Range("E88").Value = 3
$ws.Range("E89").Value = 3
$ws.Range("E90").Value = 3
$ws.Range("E91").Value = 3
$ws.Range("E92").Value = 3
$ws.Range("E93").Value = 3
$ws.Range("E94").Value = 3
$ws.Range("E95").Value = 3
$ws.Range("E96").Value = 3
$ws.Range("E97").Value = 3
